$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.025.35"
$ws.Range("D3").Value = "1.651.98"
$ws.Range("E3").Value = "  +3.54%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.75"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.884.42"
$ws.Range("E12").Value = "  +3.54%  "
$ws.Range("D13").Value = "1.646.55"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.21"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "238.99"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.40%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "27.009.44"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").Value = "1.512.53"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.04%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.579"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.17"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +10.00%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.25"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("D44").Value = "1.791.56"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.64"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0976"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.17%  "
